# Update cohort retention metrics (num_customers, cohort_size, retention_rate)
# Columns: A=cohort_year, B=period_index, C=num_customers, D=cohort_size, E=retention_rate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (C, D) new values; E is recomputed as C/D for every data row
$updates = @{
    2  = @(107, 107)
    3  = @(64,  107)
    4  = @(26,  107)
    5  = @(17,  107)
    6  = @(15,  107)
    7  = @(11,  107)
    8  = @(417, 417)
    9  = @(140, 417)
    10 = @(104, 417)
    11 = @(85,  417)
    12 = @(63,  417)
    15 = @(61,  193)
    16 = @(53,  193)
    17 = @(124, 124)
    18 = @(84,  124)
    19 = @(59,  124)
    20 = @(204, 204)
    21 = @(104, 204)
    22 = @(45,  45)
}

foreach ($row in $updates.Keys) {
    $c = $updates[$row][0]
    $d = $updates[$row][1]
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $c / $d
}
